$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativação date updated: "01/01/2020" -> "01/01/2023" ---
# Force text (not an auto-converted date serial) while keeping the cell's
# original style index: set to "@" (text) before assigning the date-like
# literal, then restore the lowercase "general" format (which this host
# resolves back to the built-in General numFmt, unlike "General" which
# would otherwise register a brand-new style entry).
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("B8").NumberFormat = "general"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2023"
$ws.Range("C8").NumberFormat = "general"

# --- New professor ("Objetivos" row 10 and "Programa resumido" row 13) ---
$ws.Range("B10").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("C10").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("B13").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("C13").Value = "5840712 - Ângelo Capri Neto"

# --- Critério text ---
$ws.Range("B19").Value = "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE."
$ws.Range("C19").Value = "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE."

# --- Norma de recuperação text ---
$ws.Range("B20").Value = "NF = (P1 + P2 + LE) /3"
$ws.Range("C20").Value = "NF = (P1 + P2 + LE) /3"

# --- Bibliografia text ---
$ws.Range("B21").Value = "Será realizada uma prova escrita valendo de zero a dez (NR) e a média final calculada pela equação: NF + NR"
$ws.Range("C21").Value = "Será realizada uma prova escrita valendo de zero a dez (NR) e a média final calculada pela equação: NF + NR"
